# Remove daft Excel reference to external file
#
# The "Changes" sheet has three formulas that reference an external
# workbook link ('[1]Level 1-0'!...). The external workbook is no longer
# needed (it points at the same "Level 1-0" sheet that already exists in
# this workbook), so re-point the formulas at the local "Level 1-0" sheet
# and break/remove the external link entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes")

# Re-point formulas that used the external workbook reference ('[1]Level 1-0')
# to the local "Level 1-0" worksheet instead.
$ws.Range("F7").Formula = "='Level 1-0'!A3"
$ws.Range("F8").Formula = "='Level 1-0'!A2"
$ws.Range("F17").Formula = "=IF(ISBLANK('Level 1-0'!A3),0,CONCATENATE(""-"",'Level 1-0'!A3))"

# Update the active selection on the sheet (matches the saved view state).
$ws.Range("F17").Select()

# Remove the now-unused external workbook link/reference from the workbook.
foreach ($link in $wb.LinkSources(1)) {
    $wb.BreakLink($link, 1)
}

$wb.Save()
